$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format so that numeric-looking strings
# (e.g. "242.81") are preserved exactly as text instead of being auto-converted
# to numbers by Excel type inference.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.879.35'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.888.32'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D5').Value = '0.7673'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').Value = '242.81'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.3124'
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('D9').Value = '25.68'
$ws.Range('D10').Value = '0.07182'
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('E11').Value = '  +4.94%  '
$ws.Range('D12').Value = '0.7641'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').Value = '1.878.74'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('D14').Value = '5.357'
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('D15').Value = '93.60'
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').Value = '6.150'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '29.943.02'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('E18').Value = '  -1.25%  '
$ws.Range('D19').Value = '244.60'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '0.000007811'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.157.43'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '0.9992'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').Value = '8.021'
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('E25').Value = '  +3.97%  '
$ws.Range('D26').Value = '9.425'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '162.87'
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').Value = '18.76'
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').Value = '2.033'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('D30').Value = '1.461'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').Value = '1.535'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = '4.503'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('D33').Value = '4.096'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').Value = '0.05450'
$ws.Range('E34').Value = '  -2.56%  '
$ws.Range('D35').Value = '1.244'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('D36').Value = '0.7432'
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('D38').Value = '2.701'
$ws.Range('E38').Value = '  +2.12%  '
$ws.Range('D39').Value = '0.01951'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('D40').Value = '2.782'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('D41').Value = '0.4463'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '1.102.03'
$ws.Range('E42').Value = '  -5.11%  '
$ws.Range('D43').Value = '73.11'
$ws.Range('E43').Value = '  -1.66%  '
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('D45').Value = '0.8524'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '102.99'
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('D48').Value = '7.681'
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('D50').Value = '3.008'
$ws.Range('E50').Value = '  -2.66%  '
$ws.Range('D51').Value = '2.050.50'
$ws.Range('E51').Value = '  +0.24%  '
